$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.850.65'
$ws.Range('E2').Value = '  +5.82%  '
$ws.Range('D3').Value = '2.224.77'
$ws.Range('E3').Value = '  +2.86%  '
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').Value = '''231.12'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.86%  '
$ws.Range('E6').Value = '  +0.43%  '
$ws.Range('D7').Value = '''60.64'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -3.11%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  +2.97%  '
$ws.Range('D10').Value = '''58.82'
$ws.Range('D10').Style = 'Normal'
$ws.Range('D11').Value = '''0.0890'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +5.15%  '
$ws.Range('E12').Value = '  -0.29%  '
$ws.Range('D13').Value = '2.553.67'
$ws.Range('E13').Value = '  +2.85%  '
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '''21.76'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.35%  '
$ws.Range('D16').Value = '''0.798'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('E17').Value = '  +1.79%  '
$ws.Range('D18').Value = '2.232.30'
$ws.Range('E18').Value = '  +2.90%  '
$ws.Range('D19').Value = '41.736.08'
$ws.Range('E19').Value = '  +5.52%  '
$ws.Range('D20').Value = '''72.78'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('D21').Value = '0.0₃0892'
$ws.Range('E21').Value = '  +0.74%  '
$ws.Range('E22').Value = '  +0.41%  '
$ws.Range('D23').Value = '''249.72'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +9.75%  '
$ws.Range('E24').Value = '  -0.05%  '
$ws.Range('E25').Value = '  +1.84%  '
$ws.Range('D26').Value = '''2.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -0.12%  '
$ws.Range('D27').Value = '''9.63'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +2.18%  '
$ws.Range('E28').Value = '  +3.22%  '
$ws.Range('D29').Value = '''167.28'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.75%  '
$ws.Range('D30').Value = '''19.92'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.52%  '
$ws.Range('E31').Value = '  -1.92%  '
$ws.Range('E32').Value = '  -2.22%  '
$ws.Range('D33').Value = '''0.121'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -0.24%  '
$ws.Range('E34').Value = '  +5.69%  '
$ws.Range('D35').Value = '''4.60'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +3.43%  '
$ws.Range('D36').Value = '''0.0625'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.29%  '
$ws.Range('D37').Value = '''6.64'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.38%  '
$ws.Range('E38').Value = '  -4.01%  '
$ws.Range('E39').Value = '  -1.43%  '
$ws.Range('D40').Value = '''0.000255'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +28.78%  '
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('E42').Value = '  +5.79%  '
$ws.Range('D43').Value = '''4.79'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -2.55%  '
$ws.Range('D44').Value = '''8.61'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +9.45%  '
$ws.Range('D45').Value = '''0.0979'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +7.21%  '
$ws.Range('E46').Value = '  +0.85%  '
$ws.Range('D47').Value = '''98.64'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.33%  '
$ws.Range('D48').Value = '1.468.25'
$ws.Range('E48').Value = '  -2.90%  '
$ws.Range('D49').Value = '''16.47'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.88%  '
$ws.Range('E50').Value = '  +0.16%  '
$ws.Range('B51').Value = 'ARBITRUM'
$ws.Range('C51').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D51').Value = '''1.07'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.39%  '
